$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.756
$ws.Range("A6").Value = -21.108
$ws.Range("A7").Value = -21.038
$ws.Range("A16").Value = -21.013
$ws.Range("A20").Value = -22.008
